$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Patients / Patient id case is not consistent across datasets / - / Updated to PATIENT_ID
$ws.Range("A14").Formula = "=A13+1"
$ws.Range("B14").Value = "Patients"
$ws.Range("C14").Value = "Patient id case is not consistent across datasets"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "Updated to PATIENT_ID"

# Row 15: Encounters / Encounter id case is not consistent across datasets / - / Updated to ENCOUNTER_ID
$ws.Range("A15").Formula = "=A14+1"
$ws.Range("B15").Value = "Encounters"
$ws.Range("C15").Value = "Encounter id case is not consistent across datasets"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "Updated to ENCOUNTER_ID"

# Resize the table to include the new rows
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:E15"))

# Update selection to match target state
$ws.Range("E16").Select()
